$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 'Última actualización: 13:54:35'
$ws.Range("A3").Value = 'Total filas: 175'
$ws.Range("C15").Value = '225_GOMEZ'
$ws.Range("C16").Value = '215A_EL PATO'
$ws.Range("A23").Value = '06:15:23'
$ws.Range("C23").Value = '16_SANTA ANA'
$ws.Range("D23").Value = 66
$ws.Range("A24").Value = '06:46:40'
$ws.Range("C24").Value = '23_HERNANDEZ'
$ws.Range("D24").Value = 35
$ws.Range("A35").Value = '07:51:40'
$ws.Range("C35").Value = '17_ROMERO'
$ws.Range("D35").Value = 9
$ws.Range("A36").Value = '07:26:49'
$ws.Range("C36").Value = '16_SANTA ANA'
$ws.Range("D36").Value = 34
$ws.Range("A43").Value = '08:14:55'
$ws.Range("C43").Value = '11_ETCHEVERRY'
$ws.Range("D43").Value = 14
$ws.Range("A44").Value = '06:58:58'
$ws.Range("C44").Value = '15_ABASTO'
$ws.Range("D44").Value = 90
$ws.Range("A45").Value = '07:51:40'
$ws.Range("C45").Value = '11_ETCHEVERRY'
$ws.Range("D45").Value = 38
$ws.Range("A46").Value = '08:14:55'
$ws.Range("C46").Value = '15_ABASTO'
$ws.Range("D46").Value = 15
$ws.Range("A66").Value = '08:57:42'
$ws.Range("C66").Value = '15X38_ABASTO'
$ws.Range("D66").Value = 21
$ws.Range("A67").Value = '08:49:06'
$ws.Range("C67").Value = '14_ABASTO'
$ws.Range("D67").Value = 29
$ws.Range("A71").Value = '08:49:06'
$ws.Range("C71").Value = '23_HERNANDEZ'
$ws.Range("D71").Value = 42
$ws.Range("A72").Value = '08:14:55'
$ws.Range("C72").Value = '16_SANTA ANA'
$ws.Range("D72").Value = 77
$ws.Range("C90").Value = '15_ABASTO'
$ws.Range("C91").Value = '14_ABASTO'
$ws.Range("A122").Value = '11:38:09'
$ws.Range("C122").Value = '15_ABASTO'
$ws.Range("D122").Value = 32
$ws.Range("A123").Value = '11:56:32'
$ws.Range("C123").Value = '16_P MOR-SANTA ANA'
$ws.Range("D123").Value = 14
$ws.Range("C137").Value = '14_ABASTO'
$ws.Range("C138").Value = '15X38_ABASTO'
$ws.Range("A151").Value = '13:28:27'
$ws.Range("C151").Value = '215A_EL PATO'
$ws.Range("D151").Value = 5
$ws.Range("A152").Value = '12:43:13'
$ws.Range("C152").Value = '14_ABASTO'
$ws.Range("D152").Value = 50
$ws.Range("A155").Value = '13:54:35'
$ws.Range("D155").Value = 0
$ws.Range("A156").Value = '13:54:35'
$ws.Range("D156").Value = 8
$ws.Range("A159").Value = '13:54:35'
$ws.Range("D159").Value = 12
$ws.Range("A160").Value = '13:54:35'
$ws.Range("D160").Value = 20
$ws.Range("A161").Value = '13:54:35'
$ws.Range("D161").Value = 22
$ws.Range("A162").Value = '13:54:35'
$ws.Range("D162").Value = 23
$ws.Range("A165").Value = '13:54:35'
$ws.Range("D165").Value = 33
$ws.Range("A166").Value = '13:54:35'
$ws.Range("D166").Value = 38
$ws.Range("A167").Value = '13:54:35'
$ws.Range("D167").Value = 40
$ws.Range("A168").Value = '13:54:35'
$ws.Range("D168").Value = 45
$ws.Range("A169").Value = '13:54:35'
$ws.Range("D169").Value = 53
$ws.Range("A170").Value = '13:54:35'
$ws.Range("C170").Value = '23_HERNANDEZ'
$ws.Range("D170").Value = 57
$ws.Range("A171").Value = '13:54:35'
$ws.Range("C171").Value = '16_SANTA ANA'
$ws.Range("D171").Value = 57
$ws.Range("A172").Value = '13:54:35'
$ws.Range("D172").Value = 60
$ws.Range("A173").Value = '13:54:35'
$ws.Range("D173").Value = 68
$ws.Range("A175").Value = '13:54:35'
$ws.Range("D175").Value = 79
$ws.Range("A176").Value = '13:54:35'
$ws.Range("B176").Value = '15:17'
$ws.Range("C176").Value = '14_ABASTO'
$ws.Range("D176").Value = 83
$ws.Range("E176").Value = 'LP1912'
$ws.Range("A177").Value = '13:54:35'
$ws.Range("B177").Value = '15:34'
$ws.Range("C177").Value = '215C_EL PATO'
$ws.Range("D177").Value = 100
$ws.Range("E177").Value = 'LP1912'
$ws.Range("A178").Value = '13:54:35'
$ws.Range("B178").Value = '15:41'
$ws.Range("C178").Value = '11_ETCHEVERRY'
$ws.Range("D178").Value = 107
$ws.Range("E178").Value = 'LP1912'
$ws.Range("A179").Value = '13:54:35'
$ws.Range("B179").Value = '15:53'
$ws.Range("C179").Value = '15X38_ABASTO'
$ws.Range("D179").Value = 119
$ws.Range("E179").Value = 'LP1912'
$ws.Range("A180").Value = '13:54:35'
$ws.Range("B180").Value = '15:53'
$ws.Range("C180").Value = '16_P MOR-SANTA ANA'
$ws.Range("D180").Value = 119
$ws.Range("E180").Value = 'LP1912'
# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)

$ws.Range("A2").Value = 'Última actualización: 13:54:35'
$ws.Range("A3").Value = 'Total filas: 29'
$ws.Range("A31").Value = '13:54:35'
$ws.Range("D31").Value = 40
$ws.Range("A32").Value = '13:54:35'
$ws.Range("D32").Value = 53
$ws.Range("A33").Value = '13:54:35'
$ws.Range("D33").Value = 60
$ws.Range("A34").Value = '13:54:35'
$ws.Range("B34").Value = '15:34'
$ws.Range("C34").Value = '215C_EL PATO'
$ws.Range("D34").Value = 100
$ws.Range("E34").Value = 'LP1912'
# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)

$ws.Range("A2").Value = 'Última actualización: 13:54:35'
$ws.Range("A3").Value = 'Total filas: 29'
$ws.Range("A19").Value = '09:42:42'
$ws.Range("C19").Value = '215A_LA PLATA'
$ws.Range("D19").Value = 48
$ws.Range("A20").Value = '08:49:06'
$ws.Range("C20").Value = '215B_LP-P MOR-1 Y 57'
$ws.Range("D20").Value = 101
$ws.Range("A30").Value = '13:54:35'
$ws.Range("D30").Value = 2
$ws.Range("A32").Value = '13:54:35'
$ws.Range("B32").Value = '14:26'
$ws.Range("D32").Value = 32
$ws.Range("B33").Value = '14:27'
$ws.Range("C33").Value = '215C_LA PLATA'
$ws.Range("D33").Value = 59
$ws.Range("E33").Value = 'L6203'
$ws.Range("A34").Value = '13:54:35'
$ws.Range("B34").Value = '15:22'
$ws.Range("C34").Value = '215A_LA PLATA'
$ws.Range("D34").Value = 88
$ws.Range("E34").Value = 'L6173'
